$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D4").Value = "2016-01-26 11:44:47"
$wsZh.Range("G4").Value = "2016-01-26 11:45:34"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D4").Value = "2016-01-26 11:44:58"
$wsDe.Range("G4").Value = "2016-01-26 11:45:54"
